$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Edité le :" timestamp
$ws.Range("B11").Value = "19/05/2016 17:26:29"

# Insert two new product rows right before the current last product row (row 18,
# "Radis"). This pushes "Radis" (and its bottom-border styling) down to row 20,
# and the two freshly inserted rows (18,19) pick up the plain row styling that
# the rest of the product rows use. The TOTAL row and the SUBTOTAL formula range,
# as well as the footer row below it, shift down automatically.
$ws.Rows("18:19").Insert()

# Row 16: Pain complet
$ws.Range("A16").Value = "Pain complet"
$ws.Range("B16").Value = "Pains"
$ws.Range("C16").Value = "A la pièce"
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 1
$ws.Range("F16").Formula = "=D16*E16"

# Row 17: Tomates grappe
$ws.Range("A17").Value = "Tomates grappe"
$ws.Range("B17").Value = "Légumes"
$ws.Range("C17").Value = "Au poids"
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 2
$ws.Range("F17").Formula = "=D17*E17"

# Row 18: Radis (now plain-styled, no longer the last row)
$ws.Range("A18").Value = "Radis"
$ws.Range("B18").Value = "Légumes"
$ws.Range("C18").Value = "A la pièce"
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 2
$ws.Range("F18").Formula = "=D18*E18"

# Row 19: Salade
$ws.Range("A19").Value = "Salade"
$ws.Range("B19").Value = "Légumes"
$ws.Range("C19").Value = "A la pièce"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 2
$ws.Range("F19").Formula = "=D19*E19"

# Row 20: Pomme de terre (now the last product row, keeps the bottom-border style)
$ws.Range("A20").Value = "Pomme de terre"
$ws.Range("B20").Value = "Légumes"
$ws.Range("C20").Value = "Au poids"
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 3
$ws.Range("F20").Formula = "=D20*E20"
